# Generate Report for Handback
#
# This script updates the localization-status workbook to reflect that the
# handback (from the localization vendor back to en-US) has completed and
# is in sync:
#   - Overview sheet: Status changes from "Ready for handoff" to
#     "Handed back: in sync with en-US" for both locale rows, and the
#     Status columns (E, F) are widened to fit the new, longer text.
#   - Per-locale sheets (zh-cn, de-de): the "Latest Target File",
#     "Latest Handback File" and "Latest Handback DateTime" columns (I, J,
#     K) are now populated for each localized file, with I turned into a
#     hyperlink (like the existing "Source File Name" column) pointing at
#     the same source .md file on GitHub. Columns C (Status), I and J are
#     widened to fit the new content.

$wb = $excel.ActiveWorkbook

$hyperlinkFontColor = 15570276   # BGR for RGB(0x64,0x95,0xED) -- matches the workbook's existing HyperLink style
$wideWidth  = 29.14              # renders as ~30 chars, matching the new Status column width
$fortyWidth = 39.166666666666664 # renders as exactly 40 chars, matching existing columns A/G

function Set-HandoffHyperlink {
    param(
        $ws,
        [string]$cellRef,
        [string]$displayText,
        [string]$targetUrl
    )

    $range = $ws.Range($cellRef)
    $range.Value = $displayText
    $ws.Hyperlinks.Add($range, $targetUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $displayText) | Out-Null
    # Hyperlinks.Add applies a generic theme-based hyperlink look; restore the
    # exact font formatting used by the workbook's own "HyperLink" style.
    $font = $range.Font
    $font.Underline = 2   # xlUnderlineStyleSingle
    $font.Color = $hyperlinkFontColor
}

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

$overview.Columns.Item(5).ColumnWidth = $wideWidth
$overview.Columns.Item(6).ColumnWidth = $wideWidth

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

Set-HandoffHyperlink $zhcn "I2" "b84f5bbe-6c67-4ae1-b690-2f3d283370b5.md" `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b8a67c6b3b985181b261f687657507234b5c8a4b/e2e/b84f5bbe-6c67-4ae1-b690-2f3d283370b5.md"
$zhcn.Range("J2").Value = "b84f5bbe-6c67-4ae1-b690-2f3d283370b5.d02e7ee50af643d5230de531e41d00ef1bca9c60.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-15 09:15:28"

Set-HandoffHyperlink $zhcn "I3" "ee2b127f-e981-41b1-82e4-fb07c0b804a4.md" `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b8a67c6b3b985181b261f687657507234b5c8a4b/e2e/ee2b127f-e981-41b1-82e4-fb07c0b804a4.md"
$zhcn.Range("J3").Value = "ee2b127f-e981-41b1-82e4-fb07c0b804a4.3f8b2a439c624253ccbc80bf855779a74d8b9e17.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-15 09:15:28"

$zhcn.Columns.Item(3).ColumnWidth = $wideWidth
$zhcn.Columns.Item(9).ColumnWidth = $fortyWidth
$zhcn.Columns.Item(10).ColumnWidth = $fortyWidth

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

Set-HandoffHyperlink $dede "I2" "b84f5bbe-6c67-4ae1-b690-2f3d283370b5.md" `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b8a67c6b3b985181b261f687657507234b5c8a4b/e2e/b84f5bbe-6c67-4ae1-b690-2f3d283370b5.md"
$dede.Range("J2").Value = "b84f5bbe-6c67-4ae1-b690-2f3d283370b5.d02e7ee50af643d5230de531e41d00ef1bca9c60.de-de.xlf"
$dede.Range("K2").Value = "2016-08-15 09:15:35"

Set-HandoffHyperlink $dede "I3" "ee2b127f-e981-41b1-82e4-fb07c0b804a4.md" `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b8a67c6b3b985181b261f687657507234b5c8a4b/e2e/ee2b127f-e981-41b1-82e4-fb07c0b804a4.md"
$dede.Range("J3").Value = "ee2b127f-e981-41b1-82e4-fb07c0b804a4.3f8b2a439c624253ccbc80bf855779a74d8b9e17.de-de.xlf"
$dede.Range("K3").Value = "2016-08-15 09:15:35"

$dede.Columns.Item(3).ColumnWidth = $wideWidth
$dede.Columns.Item(9).ColumnWidth = $fortyWidth
$dede.Columns.Item(10).ColumnWidth = $fortyWidth
